$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.8
$ws.Range("X2").Value = 17
$ws.Range("Z2").Value = 34
$ws.Range("AA2").Value = 23
$ws.Range("AB2").Value = 29
$ws.Range("AK2").Value = 19
$ws.Range("AN2").Value = 5.5
$ws.Range("AX2").Value = 11
$ws.Range("AY2").Value = 19
